$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: update the header/count values in B1:E1
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2: delete the now-invalid values in B2:E2 (Lichtwark deleted values)
$ws.Range("B2:E2").ClearContents()

# Row 3: delete B3, update C3/E3, add D3
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 30.612908304832558
$ws.Range("D3").Value = 30.483408351428295
$ws.Range("E3").Value = 16.963018195708059

# Update the selection to reflect the new data extent
$ws.Range("B1:E3").Select()
